# Dynamic values for email implemented
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Settings sheet updates ----
# B3 value changed (Volkswagen -> qwsdf)
$ws1.Range("B3").Value = "qwsdf"

# New "Description" column (C) for the car rows
$ws1.Range("C3").Value = "Make for the first selected car"
$ws1.Range("C4").Value = "Model for the first selected car"
$ws1.Range("C5").Value = "Make for the second selected car"
$ws1.Range("C6").Value = "Model for the second selected car"

# Remove the old "country" / "Germania" row entirely
$ws1.Rows("7:7").Delete() | Out-Null

# Widen description column, drop the old bestFit autosize
$ws1.Columns("C:C").ColumnWidth = 33.7109375

# Restore the (now stale) selection left behind in the source file
$ws1.Range("R12").Select() | Out-Null

# ---- New "Email" worksheet ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Email"

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Value"
$ws2.Range("C1").Value = "Description"
$ws2.Range("A1:C1").Font.Bold = $true

$ws2.Range("A2").Value = "email"
$ws2.Range("C2").Value = "Email address"

$ws2.Range("A3").Value = "name"
$ws2.Range("B3").Value = "Delia"
$ws2.Range("C3").Value = "Name of the person"

$ws2.Range("A4").Value = "subject"
$ws2.Range("B4").Value = "Assignment number 2 done"

$ws2.Range("A5").Value = "body"
$body = "Salut,`nAtasez excelul aferent punctului 2 din assignment.`nMultumesc,`nDiana Gradinaru"
$ws2.Range("B5").Value = $body
$ws2.Range("B5").WrapText = $true
$ws2.Rows(5).RowHeight = 90
$ws2.Range("C5").Value = "Body of the email"

# Hyperlink on the email address cell (also sets its displayed value/style)
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:diana.gradinaru.sincai@gmail.com", [Type]::Missing, [Type]::Missing, "diana.gradinaru.sincai@gmail.com") | Out-Null

# Column widths on the Email sheet
$ws2.Columns("B:B").ColumnWidth = 74.28515625
$ws2.Columns("C:C").ColumnWidth = 18.85546875

$ws2.Range("A3").Select() | Out-Null

# Keep "Settings" as the active/visible tab
$ws1.Activate()
